# databasdiagram.xlsx edit: add "table_language" / "table_books" (mislabeled
# "table_users" by the original author, reproduced verbatim) / "table_genre"
# blocks below the existing table_users/table_role/table_status diagram, plus
# two new "arrow" connector graphics, and tidy up a couple of leftover
# formatting quirks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-stamp F9:G13 onto the same "bordered" look used everywhere else.
#    (In the source file these five rows pointed at a slightly different
#    cellXf than the rest of the table; visually identical, but we line them
#    up with the common style so later copies are all consistent.)
# ---------------------------------------------------------------------------
$ws.Range("F6:G6").Copy() | Out-Null
$ws.Range("F9:G13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) New table_language block (columns C:D, rows 17-19)
# ---------------------------------------------------------------------------
$ws.Range("C10:D12").Copy() | Out-Null
$ws.Range("C17:D19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C17").Value = "table_language"
$ws.Range("C18").Value = "language_id"
$ws.Range("D18").Value = "INT,11,PRIMARY"
$ws.Range("C19").Value = "language_name"
$ws.Range("D19").Value = "VARCHAR,255"

# ---------------------------------------------------------------------------
# 3) New book table block (columns F:G, rows 17-29) - header row reuses the
#    bold "table title" look from F5, body rows reuse the plain bordered look
# ---------------------------------------------------------------------------
$ws.Range("F5:G5").Copy() | Out-Null
$ws.Range("F17:G17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F6:G6").Copy() | Out-Null
$ws.Range("F18:G29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F17").Value = "table_users"
$ws.Range("F18").Value = "book_id"
$ws.Range("G18").Value = " INT, 11, PRIMARY"
$ws.Range("F19").Value = "book_name"
$ws.Range("G19").Value = "VARCHAR, 255"
$ws.Range("F20").Value = "book_author"
$ws.Range("G20").Value = "VARCHAR, 255"
$ws.Range("F21").Value = "book_description"
$ws.Range("G21").Value = "text"
$ws.Range("F22").Value = "book_age_rec"
$ws.Range("G22").Value = "date"
$ws.Range("F23").Value = "language_id_fk"
$ws.Range("G23").Value = "INT, 11"
$ws.Range("F24").Value = "genre_id_fk"
$ws.Range("G24").Value = "INT, 11"
$ws.Range("F25").Value = "book_pages"
$ws.Range("G25").Value = "VARCHAR, 255"
$ws.Range("F26").Value = "book_price"
$ws.Range("G26").Value = " varchar 255"
$ws.Range("F27").Value = "book_created"
$ws.Range("G27").Value = "date"
$ws.Range("F28").Value = "book_cover"
$ws.Range("G28").Value = "VARCHAR, 255"
$ws.Range("F29").Value = "book_featured"
$ws.Range("G29").Value = "boolean"

# ---------------------------------------------------------------------------
# 4) New table_genre block (columns I:J, rows 17-19)
# ---------------------------------------------------------------------------
$ws.Range("I8:J8").Copy() | Out-Null
$ws.Range("I17:J17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I9:J9").Copy() | Out-Null
$ws.Range("I18:J19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I17").Value = "table_genre"
$ws.Range("I18").Value = "genre_id"
$ws.Range("J18").Value = "INT,11,PRIMARY"
$ws.Range("I19").Value = "genre_name"
$ws.Range("J19").Value = "VARCHAR,255"

# ---------------------------------------------------------------------------
# 5) Five leftover blank rows (30-34, columns C:D) with no border / no fill
#    (plain, unstyled look) below the new tables.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C30:D33").PasteSpecial(-4122) | Out-Null
$ws.Range("C34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Column widths - nudge the six lettered columns to their new sizes.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668   # C
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666   # D
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666   # F
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666   # G
$ws.Columns.Item(9).ColumnWidth = 13.0                 # I
$ws.Columns.Item(10).ColumnWidth = 14.833333333333334  # J

# ---------------------------------------------------------------------------
# 7) Two new "arrow" connector graphics linking the new tables, matching the
#    style of the two that already link table_users/table_role/table_status.
# ---------------------------------------------------------------------------
$pic1 = $ws.Shapes.AddPicture("xl/media/image1.png", $false, $true, 460.17, 299.8, 88.3, 44.17)
$pic1.Name = "Graphic 1"
$pic1.AlternativeText = "Arrow Right with solid fill"
$pic1.Rotation = 303.57

$pic2 = $ws.Shapes.AddPicture("xl/media/image1.png", $false, $true, 235.96, 297.66, 89.65, 43.67)
$pic2.Name = "Graphic 5"
$pic2.AlternativeText = "Arrow Right with solid fill"
$pic2.Rotation = 209.2

# ---------------------------------------------------------------------------
# 8) Selection, matching the cursor position the author left the sheet at.
# ---------------------------------------------------------------------------
$ws.Range("K24").Select() | Out-Null
